$p = $ppt.ActivePresentation

# The new slide ("DiSCoVER: top drugs (cerebellar stem cell control)") is a
# near-exact duplicate of slide 2, re-run with updated data; only one cell's
# text changed (a couple of non-breaking spaces became regular spaces).
# Recreate it by duplicating slide 2 and moving the copy to the end of the
# deck, then touch up that single cell.

$src = $p.Slides.Item(2)
$dupRange = $src.Duplicate()
$newSlide = $dupRange.Item(1)

$newSlide.MoveTo($p.Slides.Count)

$tbl = $newSlide.Shapes.Item(2).Table
$cell = $tbl.Cell(14, 4)
$cell.Shape.TextFrame.TextRange.Text = "highly selective Aurora B kinase inhibitor"
